# Update "想去人数" (interested-people count) values in column F
# across the four worksheets, per the source diff.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览";     Row = 14; Value = 1071 },
    @{ Sheet = "展览";     Row = 26; Value = 1021 },
    @{ Sheet = "展览";     Row = 28; Value = 1133 },
    @{ Sheet = "展览";     Row = 29; Value = 1064 },
    @{ Sheet = "展览";     Row = 35; Value = 1033 },
    @{ Sheet = "展览";     Row = 39; Value = 1399 },
    @{ Sheet = "展览";     Row = 40; Value = 3421 },
    @{ Sheet = "展览";     Row = 42; Value = 27 },
    @{ Sheet = "展览";     Row = 44; Value = 498 },
    @{ Sheet = "展览";     Row = 46; Value = 120 },
    @{ Sheet = "展览";     Row = 48; Value = 578 },
    @{ Sheet = "展览";     Row = 49; Value = 93 },

    @{ Sheet = "演出";     Row = 10; Value = 191 },
    @{ Sheet = "演出";     Row = 22; Value = 6847 },

    @{ Sheet = "本地生活"; Row = 4;  Value = 2037 },
    @{ Sheet = "本地生活"; Row = 9;  Value = 9031 },
    @{ Sheet = "本地生活"; Row = 10; Value = 1254 },
    @{ Sheet = "本地生活"; Row = 12; Value = 14 },

    @{ Sheet = "全部类型"; Row = 4;  Value = 2037 },
    @{ Sheet = "全部类型"; Row = 9;  Value = 1254 },
    @{ Sheet = "全部类型"; Row = 15; Value = 1071 },
    @{ Sheet = "全部类型"; Row = 26; Value = 1133 },
    @{ Sheet = "全部类型"; Row = 27; Value = 1064 },
    @{ Sheet = "全部类型"; Row = 30; Value = 191 },
    @{ Sheet = "全部类型"; Row = 32; Value = 1033 },
    @{ Sheet = "全部类型"; Row = 36; Value = 3421 },
    @{ Sheet = "全部类型"; Row = 40; Value = 498 },
    @{ Sheet = "全部类型"; Row = 42; Value = 578 },
    @{ Sheet = "全部类型"; Row = 45; Value = 93 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range("F" + $u.Row).Value = $u.Value
}
